$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 94, pushing existing rows 94..190 down to 96..192
$ws.Rows.Item(94).Insert()
$ws.Rows.Item(94).Insert()

# Fill new row 94 with new data entry
$ws.Cells.Item(94, 1).Value = 10
$ws.Cells.Item(94, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(94, 3).Value = "La Araucanía"
$ws.Cells.Item(94, 4).Value = 44587
$ws.Cells.Item(94, 5).Value = 9
$ws.Cells.Item(94, 6).Value = 100112052
$ws.Cells.Item(94, 7).Value = "Albahaca"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 40
$ws.Cells.Item(94, 11).Value = 6000
$ws.Cells.Item(94, 12).Value = 6000
$ws.Cells.Item(94, 13).Value = 6000
$ws.Cells.Item(94, 14).Value = "`$/paquete"
$ws.Cells.Item(94, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(94, 16).Value = 6000
$ws.Cells.Item(94, 17).Value = 1
$ws.Cells.Item(94, 18).Value = "Hortaliza"

# Fill new row 95 with new data entry
$ws.Cells.Item(95, 1).Value = 10
$ws.Cells.Item(95, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(95, 3).Value = "La Araucanía"
$ws.Cells.Item(95, 4).Value = 44587
$ws.Cells.Item(95, 5).Value = 9
$ws.Cells.Item(95, 6).Value = 100112052
$ws.Cells.Item(95, 7).Value = "Albahaca"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 35
$ws.Cells.Item(95, 11).Value = 5000
$ws.Cells.Item(95, 12).Value = 5000
$ws.Cells.Item(95, 13).Value = 5000
$ws.Cells.Item(95, 14).Value = "`$/paquete"
$ws.Cells.Item(95, 15).Value = "Región del Maule"
$ws.Cells.Item(95, 16).Value = 5000
$ws.Cells.Item(95, 17).Value = 1
$ws.Cells.Item(95, 18).Value = "Hortaliza"
